# Update "想去人数" (want-to-go count) figures for several conventions/events
# on the "展览" and "全部类型" worksheets, reflecting the data refresh at
# commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 131    # 苏州·无限次元夜场                     130 -> 131
$ws1.Range("F8").Value  = 148    # 苏州·漫语堂动漫嘉年华                 142 -> 148
$ws1.Range("F9").Value  = 336    # 苏州·第三届华盟国漫次元嘉年华          332 -> 336
$ws1.Range("F10").Value = 445    # 苏州·女神异闻录only同人展             444 -> 445
$ws1.Range("F11").Value = 505    # 苏州·绘时国乙2.0光夜同人only           504 -> 505
$ws1.Range("F12").Value = 136    # 张家港·META萌圆饿了                   135 -> 136
$ws1.Range("F13").Value = 11636  # 苏州·COME IN JOY 动漫品牌国潮文化节   11600 -> 11636
$ws1.Range("F14").Value = 5399   # 苏州·星部落&青铜树动漫嘉年华          5395 -> 5399

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 131    # 苏州·无限次元夜场                     130 -> 131
$ws4.Range("F10").Value = 148    # 苏州·漫语堂动漫嘉年华                 142 -> 148
$ws4.Range("F11").Value = 336    # 苏州·第三届华盟国漫次元嘉年华          332 -> 336
$ws4.Range("F12").Value = 445    # 苏州·女神异闻录only同人展             444 -> 445
$ws4.Range("F13").Value = 505    # 苏州·绘时国乙2.0光夜同人only           504 -> 505
$ws4.Range("F14").Value = 136    # 张家港·META萌圆饿了                   135 -> 136
$ws4.Range("F15").Value = 11636  # 苏州·COME IN JOY 动漫品牌国潮文化节   11600 -> 11636
$ws4.Range("F17").Value = 5399   # 苏州·星部落&青铜树动漫嘉年华          5395 -> 5399
